# This script removes the "LHE / Lahore, Pakistan" row (row 218) from the
# colo data table on Sheet1. All rows below it shift up by one, which also
# reduces the used range from A1:H332 to A1:H331 (matching the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 218 currently holds: LHE, Lahore, Pakistan, Lahore, Pakistan, PK, 31.52..., 74.40...
# Deleting the entire row shifts rows 219:332 up to 218:331.
$ws.Rows.Item(218).Delete()
